# updated filter function for empty assigned staff, CAPEX & OPEX and Link
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Requestor for row 2 changed from "Chee Ching" to "aik wei"
$ws.Range("F2").Value = "aik wei"

# Row 3 gains a "Created Date" value (8/5/2024) formatted as a short date
$ws.Range("G3").Value2 = 45509
$ws.Range("G3").NumberFormat = "mm-dd-yy"

# Update the active selection shown in the sheet view
$ws.Range("G4").Select()
